$wb = $excel.ActiveWorkbook

# --- Rename Sheet9 -> "Add Sources" ---
$ws9 = $wb.Worksheets.Item("Sheet9")
$ws9.Name = "Add Sources"

# --- Main sheet: rework the G/H column data around rows 26-28 ---
# Order matters here: "date_created" must become a shared string before the
# sheet9 additions below so the new shared-string table indices line up with
# the target workbook (date_created=218, Convert CSV=219, url=220, random gen data=221).
$wsMain = $wb.Worksheets.Item("Main")

# Shift the existing "user_id" / "INT REFERENCES users.id" pair down one row
# (row 26 -> row 27), then put the new "date_created" label on row 26.
$wsMain.Range("G27").Value = $wsMain.Range("G26").Value()
$wsMain.Range("H27").Value = $wsMain.Range("H26").Value()
$wsMain.Range("G26").Value = "date_created"
$wsMain.Range("H26").ClearContents()

# Row 28's G/H cells are no longer used at all - fully clear them (formatting
# included) so they drop out of the saved sheet entirely.
$wsMain.Range("G28").Clear()
$wsMain.Range("H28").Clear()

# --- "Add Sources" sheet: append a new source row ---
$wsAdd = $wb.Worksheets.Item("Add Sources")
$wsAdd.Range("A5").Value = "Convert CSV"
$wsAdd.Range("B5").Value = "http://www.convertcsv.com/"
$wsAdd.Range("C5").Value = "random gen data"

# --- Update selections to match the edited regions ---
$wsMain.Activate()
$wsMain.Range("G26:H27").Select() | Out-Null

$wsAdd.Activate()
$wsAdd.Range("C6").Select() | Out-Null
